$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.111.39'
$ws.Range('E2').Value = '  -0.26%  '

$ws.Range('D3').Value = '3.477.27'
$ws.Range('E3').Value = '  -0.35%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.96'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.43%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.02'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.43%  '

$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('E8').Value = '  -1.05%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.69'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.00%  '

$ws.Range('E10').Value = '  -0.72%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.385'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.05%  '

$ws.Range('D12').Value = '4.070.67'
$ws.Range('E12').Value = '  -0.32%  '

$ws.Range('E13').Value = '  -0.08%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000177'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.24%  '

$ws.Range('D15').Value = '3.479.51'
$ws.Range('E15').Value = '  -0.27%  '

$ws.Range('D16').Value = '64.099.51'
$ws.Range('E16').Value = '  -0.38%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.09'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.52%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '9.96'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.12%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.68'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.19%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.38'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.68%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '385.32'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.26%  '

$ws.Range('E22').Value = '  -0.59%  '

$ws.Range('D23').Value = '3.617.66'
$ws.Range('E23').Value = '  -0.33%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.48'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.35%  '

$ws.Range('E25').Value = '  +0.16%  '

$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.38'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -6.05%  '

$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000112'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.73%  '

$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.26%  '

$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.22'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.58%  '

$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.07'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.34%  '

$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.156'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.52%  '

$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.96'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.30%  '

$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.42'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.60%  '

$ws.Range('B34').Value = 'RenzoRestakedETH'
$ws.Range('C34').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D34').Value = '3.506.93'
$ws.Range('E34').Value = '  -0.13%  '

$ws.Range('B35').Value = 'USDe'
$ws.Range('C35').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.06%  '

$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '22.98'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.78%  '

$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.20'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.31%  '

$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.77'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.97%  '

$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '163.17'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.77%  '

$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.50'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.34%  '

$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0772'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.05%  '

$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.797'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.02%  '

$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.03%  '

$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.56'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.70%  '

$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.35'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.83%  '

$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.62'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.81%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.58'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.37%  '

$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.12'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.90%  '

$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.72'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.79%  '

$ws.Range('B50').Value = 'SuiNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.895'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.36%  '

$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.331.22'
$ws.Range('E51').Value = '  -5.05%  '
